# Update recomputed power-flow results for the 380 kV line case (rows 2-25, i.e. A=0..23)
# Columns touched: B, D, E, F, G, H, J, K, L, N (C, I, M, O stay 0; A is the index column)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 0)
$ws.Cells.Item(2, 2).Value = 1.241092252849029
$ws.Cells.Item(2, 4).Value = 0.1150248183383553
$ws.Cells.Item(2, 5).Value = 0.160876770920483
$ws.Cells.Item(2, 6).Value = 2.280691098417122
$ws.Cells.Item(2, 7).Value = 1.577955447602093
$ws.Cells.Item(2, 8).Value = 1.419242814550003
$ws.Cells.Item(2, 10).Value = 0.2394846864840634
$ws.Cells.Item(2, 11).Value = 0.6674023364164725
$ws.Cells.Item(2, 12).Value = 0.2316168997267027
$ws.Cells.Item(2, 14).Value = 2.620626096286461

# Row 3 (A3 = 1)
$ws.Cells.Item(3, 2).Value = 1.217650012112159
$ws.Cells.Item(3, 4).Value = 0.1149542754931545
$ws.Cells.Item(3, 5).Value = 0.161498686704741
$ws.Cells.Item(3, 6).Value = 2.281751497753461
$ws.Cells.Item(3, 7).Value = 1.575296100191409
$ws.Cells.Item(3, 8).Value = 1.423218395150712
$ws.Cells.Item(3, 10).Value = 0.2406192785472721
$ws.Cells.Item(3, 11).Value = 0.6128881741464056
$ws.Cells.Item(3, 12).Value = 0.2221006041640692
$ws.Cells.Item(3, 14).Value = 2.643506227855685

# Row 4 (A4 = 2)
$ws.Cells.Item(4, 2).Value = 1.20382441845797
$ws.Cells.Item(4, 4).Value = 0.1149349262874253
$ws.Cells.Item(4, 5).Value = 0.1619101887677443
$ws.Cells.Item(4, 6).Value = 2.283439278074553
$ws.Cells.Item(4, 7).Value = 1.574491615809464
$ws.Cells.Item(4, 8).Value = 1.426232704689625
$ws.Cells.Item(4, 10).Value = 0.2413564485771298
$ws.Cells.Item(4, 11).Value = 0.5796676717363596
$ws.Cells.Item(4, 12).Value = 0.2163652705041841
$ws.Cells.Item(4, 14).Value = 2.658297508788284

# Row 5 (A5 = 3)
$ws.Cells.Item(5, 2).Value = 1.198333720386586
$ws.Cells.Item(5, 4).Value = 0.1149330936478918
$ws.Cells.Item(5, 5).Value = 0.1620853511240394
$ws.Cells.Item(5, 6).Value = 2.284387863616928
$ws.Cells.Item(5, 7).Value = 1.57437200994157
$ws.Cells.Item(5, 8).Value = 1.427605309141057
$ws.Cells.Item(5, 10).Value = 0.2416670616631214
$ws.Cells.Item(5, 11).Value = 0.5661936418007372
$ws.Cells.Item(5, 12).Value = 0.2140552639017272
$ws.Cells.Item(5, 14).Value = 2.664511859329256

# Row 6 (A6 = 4)
$ws.Cells.Item(6, 2).Value = 1.197430666685705
$ws.Cells.Item(6, 4).Value = 0.1149331556607471
$ws.Cells.Item(6, 5).Value = 0.1621148885126571
$ws.Cells.Item(6, 6).Value = 2.284561130787125
$ws.Cells.Item(6, 7).Value = 1.574364723912751
$ws.Cells.Item(6, 8).Value = 1.427841944052943
$ws.Cells.Item(6, 10).Value = 0.2417192559884689
$ws.Cells.Item(6, 11).Value = 0.5639601414583524
$ws.Cells.Item(6, 12).Value = 0.2136733342526895
$ws.Cells.Item(6, 14).Value = 2.66555502891779

# Row 7 (A7 = 5)
$ws.Cells.Item(7, 2).Value = 1.203749787919179
$ws.Cells.Item(7, 4).Value = 0.1149348770307839
$ws.Cells.Item(7, 5).Value = 0.1619125207908096
$ws.Cells.Item(7, 6).Value = 2.283451014904628
$ws.Cells.Item(7, 7).Value = 1.574489159752858
$ws.Cells.Item(7, 8).Value = 1.42625063193799
$ws.Cells.Item(7, 10).Value = 0.2413605962457757
$ws.Cells.Item(7, 11).Value = 0.5794856982665806
$ws.Cells.Item(7, 12).Value = 0.216334006724253
$ws.Cells.Item(7, 14).Value = 2.658380561413885

# Row 8 (A8 = 6)
$ws.Cells.Item(8, 2).Value = 1.232891765198389
$ws.Cells.Item(8, 4).Value = 0.1149955369103957
$ws.Cells.Item(8, 5).Value = 0.1610850657545324
$ws.Cells.Item(8, 6).Value = 2.280841623429083
$ws.Cells.Item(8, 7).Value = 1.576866529025338
$ws.Cells.Item(8, 8).Value = 1.420494661468879
$ws.Cells.Item(8, 10).Value = 0.2398674946083839
$ws.Cells.Item(8, 11).Value = 0.6485538787785288
$ws.Cells.Item(8, 12).Value = 0.2283133844255474
$ws.Cells.Item(8, 14).Value = 2.628361004115959

# Row 9 (A9 = 7)
$ws.Cells.Item(9, 2).Value = 1.294528229009558
$ws.Cells.Item(9, 4).Value = 0.1153035529592579
$ws.Cells.Item(9, 5).Value = 0.1596968441071196
$ws.Cells.Item(9, 6).Value = 2.283947648192608
$ws.Cells.Item(9, 7).Value = 1.588106867952732
$ws.Cells.Item(9, 8).Value = 1.413752960273328
$ws.Cells.Item(9, 10).Value = 0.2372601726892345
$ws.Cells.Item(9, 11).Value = 0.7859821657690418
$ws.Cells.Item(9, 12).Value = 0.2526564273957774
$ws.Cells.Item(9, 14).Value = 2.575385440001725

# Row 10 (A10 = 8)
$ws.Cells.Item(10, 2).Value = 1.342531341851014
$ws.Cells.Item(10, 4).Value = 0.11564375505003
$ws.Cells.Item(10, 5).Value = 0.1588187528059368
$ws.Cells.Item(10, 6).Value = 2.291242781475347
$ws.Cells.Item(10, 7).Value = 1.600386711131421
$ws.Cells.Item(10, 8).Value = 1.411568397081311
$ws.Cells.Item(10, 10).Value = 0.2355387635080586
$ws.Cells.Item(10, 11).Value = 0.8881600816483797
$ws.Cells.Item(10, 12).Value = 0.2710583259779327
$ws.Cells.Item(10, 14).Value = 2.540054075858787

# Row 11 (A11 = 9)
$ws.Cells.Item(11, 2).Value = 1.364956183342997
$ws.Cells.Item(11, 4).Value = 0.1158229771540746
$ws.Cells.Item(11, 5).Value = 0.1584498591940573
$ws.Cells.Item(11, 6).Value = 2.295650412534485
$ws.Cells.Item(11, 7).Value = 1.606849047635848
$ws.Cells.Item(11, 8).Value = 1.411175323675479
$ws.Cells.Item(11, 10).Value = 0.2347975432223368
$ws.Cells.Item(11, 11).Value = 0.9349064935367437
$ws.Cells.Item(11, 12).Value = 0.2795418087394665
$ws.Cells.Item(11, 14).Value = 2.524759880147045

# Row 12 (A12 = 10)
$ws.Cells.Item(12, 2).Value = 1.373532023526764
$ws.Cells.Item(12, 4).Value = 0.1158943370894576
$ws.Cells.Item(12, 5).Value = 0.1583145447546581
$ws.Cells.Item(12, 6).Value = 2.297476035043474
$ws.Cells.Item(12, 7).Value = 1.609422314324348
$ws.Cells.Item(12, 8).Value = 1.411112801820565
$ws.Cells.Item(12, 10).Value = 0.2345228612473882
$ws.Cells.Item(12, 11).Value = 0.9526460603857458
$ws.Cells.Item(12, 12).Value = 0.2827703695230923
$ws.Cells.Item(12, 14).Value = 2.519080282820127

# Row 13 (A13 = 11)
$ws.Cells.Item(13, 2).Value = 1.371681333366752
$ws.Cells.Item(13, 4).Value = 0.1158788135242119
$ws.Cells.Item(13, 5).Value = 0.1583434927026364
$ws.Cells.Item(13, 6).Value = 2.297075892885786
$ws.Cells.Item(13, 7).Value = 1.608862504034562
$ws.Cells.Item(13, 8).Value = 1.411122428448522
$ws.Cells.Item(13, 10).Value = 0.2345817522510236
$ws.Cells.Item(13, 11).Value = 0.948823856471563
$ws.Cells.Item(13, 12).Value = 0.2820743290842671
$ws.Cells.Item(13, 14).Value = 2.520298504020353

# Row 14 (A14 = 12)
$ws.Cells.Item(14, 2).Value = 1.365660041529679
$ws.Cells.Item(14, 4).Value = 0.1158287781296821
$ws.Cells.Item(14, 5).Value = 0.1584386391527381
$ws.Cells.Item(14, 6).Value = 2.295797470439979
$ws.Cells.Item(14, 7).Value = 1.607058223466737
$ws.Cells.Item(14, 8).Value = 1.41116845015037
$ws.Cells.Item(14, 10).Value = 0.2347748247779959
$ws.Cells.Item(14, 11).Value = 0.9363651848844654
$ws.Cells.Item(14, 12).Value = 0.2798071033313079
$ws.Cells.Item(14, 14).Value = 2.524290371511142

# Row 15 (A15 = 13)
$ws.Cells.Item(15, 2).Value = 1.361982755651866
$ws.Cells.Item(15, 4).Value = 0.1157985840752502
$ws.Cells.Item(15, 5).Value = 0.1584974887119541
$ws.Cells.Item(15, 6).Value = 2.295034785410351
$ws.Cells.Item(15, 7).Value = 1.60596947788099
$ws.Cells.Item(15, 8).Value = 1.411207880491133
$ws.Cells.Item(15, 10).Value = 0.2348938684186663
$ws.Cells.Item(15, 11).Value = 0.9287387881239511
$ws.Cells.Item(15, 12).Value = 0.2784204490688182
$ws.Cells.Item(15, 14).Value = 2.526750091368804

# Row 16 (A16 = 14)
$ws.Cells.Item(16, 2).Value = 1.341077613672098
$ws.Cells.Item(16, 4).Value = 0.1156325325479628
$ws.Cells.Item(16, 5).Value = 0.1588434743239162
$ws.Cells.Item(16, 6).Value = 2.290976641989246
$ws.Cells.Item(16, 7).Value = 1.599982023963676
$ws.Cells.Item(16, 8).Value = 1.411606168335453
$ws.Cells.Item(16, 10).Value = 0.2355880448195471
$ws.Cells.Item(16, 11).Value = 0.8851103939804545
$ws.Cells.Item(16, 12).Value = 0.2705061622747138
$ws.Cells.Item(16, 14).Value = 2.541069253661249

# Row 17 (A17 = 15)
$ws.Cells.Item(17, 2).Value = 1.328403218396147
$ws.Cells.Item(17, 4).Value = 0.115536912695088
$ws.Cells.Item(17, 5).Value = 0.1590635395162545
$ws.Cells.Item(17, 6).Value = 2.288765992231973
$ws.Cells.Item(17, 7).Value = 1.596533414370114
$ws.Cells.Item(17, 8).Value = 1.412004319519241
$ws.Cells.Item(17, 10).Value = 0.2360246081091404
$ws.Cells.Item(17, 11).Value = 0.8584134067600928
$ws.Cells.Item(17, 12).Value = 0.265679706367024
$ws.Cells.Item(17, 14).Value = 2.550052969206853

# Row 18 (A18 = 16)
$ws.Cells.Item(18, 2).Value = 1.321168615888439
$ws.Cells.Item(18, 4).Value = 0.1154842187065839
$ws.Cells.Item(18, 5).Value = 0.1591929920894994
$ws.Cells.Item(18, 6).Value = 2.287596993507577
$ws.Cells.Item(18, 7).Value = 1.594632327310777
$ws.Cells.Item(18, 8).Value = 1.412289868882283
$ws.Cells.Item(18, 10).Value = 0.2362796485890124
$ws.Cells.Item(18, 11).Value = 0.8430829892155884
$ws.Cells.Item(18, 12).Value = 0.2629142395257276
$ws.Cells.Item(18, 14).Value = 2.555293409654027

# Row 19 (A19 = 17)
$ws.Cells.Item(19, 2).Value = 1.318728630025106
$ws.Cells.Item(19, 4).Value = 0.1154667739170954
$ws.Cells.Item(19, 5).Value = 0.1592373171093699
$ws.Cells.Item(19, 6).Value = 2.287218798252454
$ws.Cells.Item(19, 7).Value = 1.594002810966018
$ws.Cells.Item(19, 8).Value = 1.412396263999241
$ws.Cells.Item(19, 10).Value = 0.2363666782726739
$ws.Cells.Item(19, 11).Value = 0.8378966768391933
$ws.Cells.Item(19, 12).Value = 0.2619797209043355
$ws.Cells.Item(19, 14).Value = 2.55708031038677

# Row 20 (A20 = 18)
$ws.Cells.Item(20, 2).Value = 1.329746702078069
$ws.Cells.Item(20, 4).Value = 0.1155468533182358
$ws.Cells.Item(20, 5).Value = 0.1590398155706625
$ws.Cells.Item(20, 6).Value = 2.28899071081527
$ws.Cells.Item(20, 7).Value = 1.596891989322359
$ws.Cells.Item(20, 8).Value = 1.411956084005652
$ws.Cells.Item(20, 10).Value = 0.2359777274651504
$ws.Cells.Item(20, 11).Value = 0.8612527639496932
$ws.Cells.Item(20, 12).Value = 0.2661923963291741
$ws.Cells.Item(20, 14).Value = 2.549089056236944

# Row 21 (A21 = 19)
$ws.Cells.Item(21, 2).Value = 1.367426364170313
$ws.Cells.Item(21, 4).Value = 0.1158433801550487
$ws.Cells.Item(21, 5).Value = 0.158410573656127
$ws.Cells.Item(21, 6).Value = 2.29616872584819
$ws.Cells.Item(21, 7).Value = 1.607584760521348
$ws.Cells.Item(21, 8).Value = 1.411152589996817
$ws.Cells.Item(21, 10).Value = 0.2347179519368767
$ws.Cells.Item(21, 11).Value = 0.9400235784574704
$ws.Cells.Item(21, 12).Value = 0.2804726077622917
$ws.Cells.Item(21, 14).Value = 2.523114823465505

# Row 22 (A22 = 20)
$ws.Cells.Item(22, 2).Value = 1.392541810716153
$ws.Cells.Item(22, 4).Value = 0.1160575224798492
$ws.Cells.Item(22, 5).Value = 0.1580248360861694
$ws.Cells.Item(22, 6).Value = 2.301772433015387
$ws.Cells.Item(22, 7).Value = 1.615308289093008
$ws.Cells.Item(22, 8).Value = 1.411130613782859
$ws.Cells.Item(22, 10).Value = 0.2339295926267351
$ws.Cells.Item(22, 11).Value = 0.9917245003211406
$ws.Cells.Item(22, 12).Value = 0.289899053056871
$ws.Cells.Item(22, 14).Value = 2.506791919227148

# Row 23 (A23 = 21)
$ws.Cells.Item(23, 2).Value = 1.379092593811777
$ws.Cells.Item(23, 4).Value = 0.1159413774189488
$ws.Cells.Item(23, 5).Value = 0.1582283827675104
$ws.Cells.Item(23, 6).Value = 2.298698156355712
$ws.Cells.Item(23, 7).Value = 1.611118784726642
$ws.Cells.Item(23, 8).Value = 1.411096321362862
$ws.Cells.Item(23, 10).Value = 0.2343471600185838
$ws.Cells.Item(23, 11).Value = 0.9641107851060156
$ws.Cells.Item(23, 12).Value = 0.2848594620488285
$ws.Cells.Item(23, 14).Value = 2.515444007020285

# Row 24 (A24 = 22)
$ws.Cells.Item(24, 2).Value = 1.329139150853223
$ws.Cells.Item(24, 4).Value = 0.1155423520626009
$ws.Cells.Item(24, 5).Value = 0.1590505320262929
$ws.Cells.Item(24, 6).Value = 2.288888798016899
$ws.Cells.Item(24, 7).Value = 1.596729623522648
$ws.Cells.Item(24, 8).Value = 1.41197771481832
$ws.Cells.Item(24, 10).Value = 0.2359989095744783
$ws.Cells.Item(24, 11).Value = 0.8599690343554869
$ws.Cells.Item(24, 12).Value = 0.2659605801377012
$ws.Cells.Item(24, 14).Value = 2.54952460584882

# Row 25 (A25 = 23)
$ws.Cells.Item(25, 2).Value = 1.277375195084232
$ws.Cells.Item(25, 4).Value = 0.1152001269130629
$ws.Cells.Item(25, 5).Value = 0.160047410206519
$ws.Cells.Item(25, 6).Value = 2.282227045452416
$ws.Cells.Item(25, 7).Value = 1.584360646745665
$ws.Cells.Item(25, 8).Value = 1.41509039526882
$ws.Cells.Item(25, 10).Value = 0.2379313282830404
$ws.Cells.Item(25, 11).Value = 0.7485917125013657
$ws.Cells.Item(25, 12).Value = 0.2459800412287052
$ws.Cells.Item(25, 14).Value = 2.589086159375437
